$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = 44187
$districts = @(
    @('Bánovce nad Bebravou', 41),
    @('Banská Bystrica', 88),
    @('Banská Štiavnica', 5),
    @('Bardejov', 13),
    @('Bratislava', 259),
    @('Brezno', 30),
    @('Bytča', 17),
    @('Čadca', 7),
    @('Detva', 12),
    @('Dolný Kubín', 28),
    @('Dunajská Streda', 26),
    @('Galanta', 19),
    @('Gelnica', 11),
    @('Hlohovec', 27),
    @('Humenné', 25),
    @('Ilava', 74),
    @('Kežmarok', 37),
    @('Komárno', 44),
    @('Košice', 182),
    @('Košice - okolie', 86),
    @('Krupina', 7),
    @('Kysucké Nové Mesto', 14),
    @('Levice', 46),
    @('Levoča', 18),
    @('Liptovský Mikuláš', 164),
    @('Lučenec', 45),
    @('Malacky', 26),
    @('Martin', 103),
    @('Medzilaborce', 1),
    @('Michalovce', 33),
    @('Myjava', 25),
    @('Námestovo', 2),
    @('Nitra', 480),
    @('Nové Mesto nad Váhom', 136),
    @('Nové Zámky', 56),
    @('Partizánske', 12),
    @('Pezinok', 18),
    @('Piešťany', 41),
    @('Poltár', 14),
    @('Poprad', 79),
    @('Považská Bystrica', 129),
    @('Prešov', 158),
    @('Prievidza', 87),
    @('Púchov', 84),
    @('Revúca', 8),
    @('Rimavská Sobota', 10),
    @('Rožňava', 12),
    @('Ružomberok', 51),
    @('Sabinov', 18),
    @('Senec', 24),
    @('Senica', 60),
    @('Skalica', 46),
    @('Snina', 3),
    @('Sobrance', 7),
    @('Spišská Nová Ves', 73),
    @('Stará Ľubovňa', 4),
    @('Stropkov', 10),
    @('Svidník', 34),
    @('Šaľa', 27),
    @('Topoľčany', 7),
    @('Trebišov', 34),
    @('Trenčín', 85),
    @('Trnava', 92),
    @('Turčianske Teplice', 4),
    @('Veľký Krtíš', 8),
    @('Vranov nad Topľou', 39),
    @('Zlaté Moravce', 13),
    @('Zvolen', 55),
    @('Žarnovica', 3),
    @('Žiar nad Hronom', 22),
    @('Žilina', 129)
)

$startRow = 6639
for ($i = 0; $i -lt $districts.Length; $i++) {
    $row = $startRow + $i
    $name = $districts[$i][0]
    $count = $districts[$i][1]
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $count
}
